# Updated cryptos list on Thu Oct  3 07:56:18 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.970.39"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "2.366.57"
$ws.Range("E3").Value = "  -4.63%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.51"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.28"
$ws.Range("E6").Value = "  -5.23%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  -13.04%  "
$ws.Range("D9").Value = "2.360.00"
$ws.Range("E9").Value = "  -4.81%  "
$ws.Range("E10").Value = "  -3.57%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  -4.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.343"
$ws.Range("E13").Value = "  -4.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.97"
$ws.Range("E14").Value = "  -6.18%  "
$ws.Range("D15").Value = "2.798.03"
$ws.Range("E15").Value = "  -4.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000163"
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("D17").Value = "60.902.92"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "2.371.82"
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.64"
$ws.Range("E19").Value = "  -5.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.09"
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.73"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.60"
$ws.Range("E22").Value = "  -8.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.87"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.47"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.08"
$ws.Range("E26").Value = "  +3.46%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "2.488.96"
$ws.Range("E28").Value = "  -4.53%  "
$ws.Range("D29").Value = "0.0₃0914"
$ws.Range("E29").Value = "  -9.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "519.62"
$ws.Range("E30").Value = "  -8.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.01"
$ws.Range("E31").Value = "  -4.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.145"
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("E34").Value = "  -5.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.61"
$ws.Range("E37").Value = "  -7.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.36"
$ws.Range("E38").Value = "  -10.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.372"
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("E40").Value = "  +3.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.09"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "139.95"
$ws.Range("E42").Value = "  -3.29%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.28"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "139.32"
$ws.Range("E45").Value = "  -6.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.11"
$ws.Range("E46").Value = "  -14.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.57"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0513"
$ws.Range("E48").Value = "  -5.69%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.74"
$ws.Range("E49").Value = "  -10.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.570"
$ws.Range("E50").Value = "  -5.18%  "
